{"js": "// The document's title paragraph reads \"Online Calculator\" and should be\n// shortened to just \"Calculator\" (the rest of the document body is\n// unchanged). Find the exact run of text and replace it in place so the\n// surrounding run formatting (bold, font size 28, centered heading, etc.)\n// is preserved.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"Online Calculator\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  // Replace just the title occurrence; this keeps the existing run's\n  // character formatting (bold/size/color) intact.\n  searchResults.items[0].insertText(\"Calculator\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document's title paragraph reads \"Online Calculator\" and should be\n# shortened to just \"Calculator\" (the rest of the document body is\n# unchanged). Use Find/Replace scoped to the document content so the\n# existing run formatting (bold, size 28, centered heading, etc.) is kept.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Online Calculator\"\n$find.Replacement.Text = \"Calculator\"\n$find.Execute(\n    \"Online Calculator\",  # FindText\n    $true,                  # MatchCase\n    $false,                # MatchWholeWord\n    $false,                # MatchWildcards\n    $false,                # MatchSoundsLike\n    $false,                # MatchAllWordForms\n    $true,                 # Forward\n    1,                      # Wrap (wdFindContinue)\n    $false,                # Format\n    \"Calculator\",           # ReplaceWith\n    1                       # Replace (wdReplaceOne)\n)\n"}
